# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# header style from column G and filling the data rows with 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: text "Save", formatted like the other headers (copy G1's format).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells H2:H9: numeric 0, no special style (matches columns B-G's data cells).
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
